$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.162357926368713
$ws.Range("B1").Value = 2.417994737625122
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.379735231399536
$ws.Range("E1").Value = 1.233781933784485
